$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the country labels for Honduras / Portugal (rows 53 & 54) ---
# Row 53 currently shows "Portugal" -> becomes "Honduras"
# Row 54 currently shows "Honduras" -> becomes "Portugal"
$ws.Range("A53").Value = "Honduras"
$ws.Range("A54").Value = "Portugal"

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 05:34"

# --- Update numeric data for the affected rows ---

# Row 24 (Pakistan)
$ws.Range("B24").Value = 315727
$ws.Range("C24").Value = 467
$ws.Range("D24").Value = 300616
$ws.Range("E24").Value = 8588
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 6523

# Row 35 (Belgica)
$ws.Range("B35").Value = 132203
$ws.Range("C35").Value = 1968
$ws.Range("D35").Value = 19712
$ws.Range("E35").Value = 102413
$ws.Range("G35").Value = 14
$ws.Range("H35").Value = 10078

# Row 39 (Kazajistan)
$ws.Range("B39").Value = 108296
$ws.Range("C39").Value = 60
$ws.Range("D39").Value = 103367
$ws.Range("E39").Value = 3204

# Row 53 (now Honduras)
$ws.Range("B53").Value = 80020
$ws.Range("C53").Value = 391
$ws.Range("D53").Value = 29768
$ws.Range("E53").Value = 47819
$ws.Range("G53").Value = 11
$ws.Range("H53").Value = 2433

# Row 54 (now Portugal)
$ws.Range("B54").Value = 79885
$ws.Range("D54").Value = 50454
$ws.Range("E54").Value = 27413
$ws.Range("H54").Value = 2018

# Row 81 (Hungria)
$ws.Range("B81").Value = 27173
$ws.Range("C81").Value = 24
$ws.Range("E81").Value = 1386
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 895

# Row 133 (Guinea Ecuatorial)
$ws.Range("B133").Value = 4767
$ws.Range("H133").Value = 82

# Row 173 (Santo Tome y Principe)
$ws.Range("B173").Value = 686
$ws.Range("C173").Value = 7
$ws.Range("D173").Value = 592
$ws.Range("E173").Value = 72

# Row 186 (San Martin (Parte Francesa))
$ws.Range("B186").Value = 315
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 8
